$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2: First name (Nombre) additional passengers ---
$ws.Range("C2").Value = "Nombre p2"
$ws.Range("D2").Value = "Nombre p3"
$ws.Range("E2").Value = "Nombre p4"

# --- Row 3: Last name (Apellido) additional passengers ---
$ws.Range("C3").Value = "Apellido p2"
$ws.Range("D3").Value = "Apellido p3"
$ws.Range("E3").Value = "Apellido p4"

# --- Row 4: Meal options for additional passengers ---
$ws.Range("C4").Value = "Bland"
$ws.Range("D4").Value = "Low Calorie"
$ws.Range("E4").Value = "Diabetic"

# Match formatting of the existing columns for the newly added cells
$ws.Range("B2").Copy()
$ws.Range("C2:E3").PasteSpecial(-4122)

$ws.Range("B4").Copy()
$ws.Range("C4:E4").PasteSpecial(-4122)

$excel.CutCopyMode = 0

# --- Column C is narrower now that more columns carry data ---
$ws.Columns.Item(3).ColumnWidth = 10.763

# --- Sheet dimension / view bookkeeping ---
$ws.Range("A1").Select()
$ws.Range("E4").Select()
